$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 238, shifting existing rows 238-292 down to 239-293
$ws.Rows.Item(238).Insert()

# Populate the newly inserted row with the new record
$ws.Cells.Item(238, 1).Value = 11
$ws.Cells.Item(238, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(238, 3).Value = "Bíobío"
$ws.Cells.Item(238, 4).Value = 45258
$ws.Cells.Item(238, 5).Value = 8
$ws.Cells.Item(238, 6).Value = "Fruta"
$ws.Cells.Item(238, 7).Value = 100102
$ws.Cells.Item(238, 8).Value = "Cítricos"
$ws.Cells.Item(238, 9).Value = 100102004
$ws.Cells.Item(238, 10).Value = "Mandarina"
$ws.Cells.Item(238, 11).Value = "Murcott"
$ws.Cells.Item(238, 12).Value = "Primera"
$ws.Cells.Item(238, 13).Value = 100
$ws.Cells.Item(238, 14).Value = 9000
$ws.Cells.Item(238, 15).Value = 10000
$ws.Cells.Item(238, 16).Value = 9500
$ws.Cells.Item(238, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(238, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(238, 19).Value = 528
$ws.Cells.Item(238, 20).Value = 18
